$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the OSVPPRU10 -> OSVPPRU04 code references
$ws.Range("G2").Value = "OSVPPRU04"
$ws.Range("G3").Value = "OSVPPRU04"

# Remove trailing period from the invalid user/password message
$ws.Range("E5").Value = "Usuario o clave inválida. Inténtalo nuevamente"
$ws.Range("E6").Value = "Usuario o clave inválida. Inténtalo nuevamente"
$ws.Range("E8").Value = "Usuario o clave inválida. Inténtalo nuevamente"

# Adjust the sheet view: scroll so column B is at top-left, and select C12
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C12").Select()
